$d = $word.ActiveDocument
$full = $d.Content.Text
$marker = " inserita nella lista delle prenotazioni all’interno del campo "
$endmarker = ", prezzo, tempo e scali"
$idxStart = $full.IndexOf($marker)
$idxMarkerEnd = $idxStart + $marker.Length
$idxEndMarker = $full.IndexOf($endmarker, $idxStart)
$idxEnd = $idxEndMarker + $endmarker.Length

# Capture formatting reference BEFORE editing (1-char probe range within the "good" run)
$probe = $d.Range($idxStart, $idxStart + 1)

$rng2 = $d.Range($idxMarkerEnd, $idxEnd)
$rng2.Delete()
$ins = $d.Range($idxMarkerEnd, $idxMarkerEnd)
$ins.InsertAfter("t_utente del nodo utenteCorrente, le cui informazioni andranno inserite nei vari campi della struttura t_prenotazione, prezzo, tempo e scali")
$newRng = $d.Range($idxMarkerEnd, $idxMarkerEnd + 200)
$newRng.FormattedText = $probe.FormattedText
Write-Output "Done"
